$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14, pushing existing rows 14-23 down to 15-24.
# Excel's Insert() copies formatting from the row above, matching the
# original file's row-14 style (date column D keeps style index 2).
$ws.Rows.Item(14).Insert()

# Populate the new row 14 with the weekly record (same market/category/etc.
# as the rest of the sheet, with updated date, volume, prices).
$ws.Range("A14").Value = 1
$ws.Range("B14").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C14").Value = "Arica y Parinacota"
$ws.Range("D14").Value = 44489
$ws.Range("E14").Value = 15
$ws.Range("F14").Value = 100112013
$ws.Range("G14").Value = "Alcachofa"
$ws.Range("H14").Value = "Madrigal"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 100
$ws.Range("K14").Value = 13000
$ws.Range("L14").Value = 14000
$ws.Range("M14").Value = 13500
$ws.Range("N14").Value = "$/caja 40 unidades"
$ws.Range("O14").Value = "Región de Coquimbo"
$ws.Range("P14").Value = 338
$ws.Range("Q14").Value = 40
$ws.Range("R14").Value = "Hortaliza"
